# Update forecast score values with refreshed model output
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecasts")

$ws.Range("B3").Value = 0.6459373893103791
$ws.Range("E3").Value = 0.8380898705323252

$ws.Range("B4").Value = 0.09960389109914584
$ws.Range("E4").Value = 0.9987932693766906

$ws.Range("B5").Value = 6.433787738173758
$ws.Range("E5").Value = 83.70785218204684

$ws.Range("B6").Value = 0.9870765
$ws.Range("C6").Value = 0.0129235
$ws.Range("E6").Value = 0.4795524
$ws.Range("F6").Value = 0.5204476

$ws.Range("B7").Value = 33.7868016
$ws.Range("C7").Value = 16.1652718
$ws.Range("E7").Value = 24.0421396
$ws.Range("F7").Value = 24.8297674
